$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Alpha/Beta "Dose (pM)" readings for D11 and D12 as read off the
# experimental time course.
$ws.Range("D11").Value = 6120
$ws.Range("D12").Value = 625

# Update the active selection to match the last-edited cell.
$ws.Range("D12").Select()
